$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "student"
$ws.Range("B3").Value = "student"
$ws.Range("C5").Select() | Out-Null
